# Auto-generated edit script: refresh the crypto price/volume snapshot
# (columns D "Price" and E "Volume(1h)") for rows 2-51 on Sheet1.
#
# Values in column D are stored as plain text (not numbers) in the workbook,
# since prices use "." as a thousands separator (e.g. "42.037.00") which is
# not valid numeric syntax. For the handful of values that *are* valid
# numeric literals (e.g. "311.55", "1.00", "0.0910"), a leading apostrophe
# is used -- exactly like typing into Excel by hand -- so Excel keeps them
# as literal text (preserving trailing zeros) instead of silently
# reinterpreting them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "42.037.00"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "2.303.16"
$ws.Range("E3").Value = "  -2.50%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'311.55"
$ws.Range("E5").Value = "  -5.81%  "
$ws.Range("D6").Value = "'105.88"
$ws.Range("E6").Value = "  +5.25%  "
$ws.Range("D7").Value = "'0.621"
$ws.Range("E7").Value = "  -2.51%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -4.03%  "
$ws.Range("D10").Value = "'40.13"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").Value = "'0.0910"
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("E12").Value = "  -2.39%  "
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Value = "'0.980"
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("D15").Value = "'15.49"
$ws.Range("E15").Value = "  -5.29%  "
$ws.Range("D16").Value = "2.646.37"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").Value = "2.300.95"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").Value = "42.037.31"
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("D19").Value = "'7.65"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("D21").Value = "'74.58"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("E22").Value = "  -6.77%  "
$ws.Range("D23").Value = "'258.72"
$ws.Range("E23").Value = "  -3.97%  "
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").Value = "'9.26"
$ws.Range("E25").Value = "  -6.88%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "'10.98"
$ws.Range("E27").Value = "  -4.37%  "
$ws.Range("E28").Value = "  +3.28%  "
$ws.Range("D29").Value = "'22.75"
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("D30").Value = "'35.45"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").Value = "'0.0889"
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").Value = "'162.35"
$ws.Range("E32").Value = "  -7.78%  "
$ws.Range("D33").Value = "'2.91"
$ws.Range("E33").Value = "  -5.37%  "
$ws.Range("E34").Value = "  -3.67%  "
$ws.Range("E35").Value = "  -2.93%  "
$ws.Range("D36").Value = "'0.116"
$ws.Range("E36").Value = "  +10.27%  "
$ws.Range("D37").Value = "'4.53"
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("D38").Value = "'0.0352"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").Value = "'2.73"
$ws.Range("E39").Value = "  -6.66%  "
$ws.Range("D40").Value = "'3.62"
$ws.Range("E40").Value = "  -5.53%  "
$ws.Range("D41").Value = "'98.13"
$ws.Range("E41").Value = "  +8.52%  "
$ws.Range("D42").Value = "'1.46"
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("D43").Value = "'70.18"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "'0.230"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "'12.05"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("D47").Value = "'111.29"
$ws.Range("E47").Value = "  -5.68%  "
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("D49").Value = "'8.91"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("D50").Value = "'73.60"
$ws.Range("E50").Value = "  +4.69%  "
$ws.Range("E51").Value = "  -0.26%  "
